$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16 (A16 = "cuando se anula una OT consultar si se quieren anular la OT asociadas")
# changes status from "no comenzado" to "terminado"
$ws.Range("B16").Value = "terminado"

# New rows 18-20 with new tasks, status "terminado"
$ws.Range("A18").Value = "modificar en reportes de produccion filtrar ordenes anuladas"
$ws.Range("B18").Value = "terminado"

$ws.Range("A19").Value = "reportes de ventas, fcs en dolares multiplicar por TC"
$ws.Range("B19").Value = "terminado"

$ws.Range("A20").Value = "permitir modificaciond e cantidades en ots"
$ws.Range("B20").Value = "terminado"

# Update view: scroll position and active selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A21").Select()
